$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the count values for existing rows
$ws.Range("B2").Value = 26880
$ws.Range("B3").Value = 2735
$ws.Range("B4").Value = 21741
$ws.Range("B5").Value = 11008

# Remove row 6 (label 4, count 629) entirely
$ws.Rows.Item(6).Delete()
